# Diagrama de componentes versao final (atualizado)
# Applies the geometry / text-size tweaks described by the commit diff.
#
# NOTE on numeric literals: PowerPoint's Shape geometry properties
# (Left/Top/Width/Height) are single-precision (float32) point values
# that get converted to EMU (1 pt = 12700 EMU) by truncation. The
# literals below were back-solved so that, after that float32 ->
# truncate-to-EMU pipeline, they land exactly on the EMU values from
# the target OOXML.

function Get-ShapeById {
    param($Shapes, [int]$Id)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $cand = $Shapes.Item($i)
        if ($cand.Id -eq $Id) {
            return $cand
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape Id=3 "Retângulo 2" (o bloco "Banco de Dados") ---------------
$shape3 = Get-ShapeById $s.Shapes 3
$shape3.Left   = 238.5842514038086
$shape3.Top    = 20.332362174987793
$shape3.Width  = 115.71157455444336
$shape3.Height = 90.36771774291992

# Texto do rodapé do mesmo retângulo: 11pt -> 10pt
$tr3 = $shape3.TextFrame.TextRange
for ($i = 1; $i -le $tr3.Paragraphs().Count; $i++) {
    $para = $tr3.Paragraphs($i)
    if ($para.Text.TrimEnd("`r") -eq "Componente que gerencia as conexões e transações com o Banco de Dados") {
        $para.Font.Size = 10
    }
}

# --- Connector Id=25 "Conector: Angulado 24" ----------------------------
$conn25 = Get-ShapeById $s.Shapes 25
$conn25.Left   = 248.56890106201172
$conn25.Top    = 133.15976715087893
$conn25.Width  = 72.23850631713867
$conn25.Height = 29.073227882385257
$conn25.Adjustments.Item(1) = 0.33757

# --- Connector Id=30 "Conector: Angulado 29" ----------------------------
$conn30 = Get-ShapeById $s.Shapes 30
$conn30.Left   = 118.57716751098634
$conn30.Top    = 65.51613998413086
$conn30.Width  = 120.00708389282228
$conn30.Height = 0.11929133906960489

# --- Shape Id=24 "Retângulo 23" -----------------------------------------
$shape24 = Get-ShapeById $s.Shapes 24
$shape24.Left   = 479.5991363525391
$shape24.Top    = 183.55677032470706
$shape24.Width  = 109.48677444458009
$shape24.Height = 86.10512161254884

# --- Shape Id=35 "Retângulo 34" -----------------------------------------
$shape35 = Get-ShapeById $s.Shapes 35
$shape35.Left   = 361.33967590332037
$shape35.Top    = 183.55677032470706
$shape35.Width  = 115.27086257934572
$shape35.Height = 86.55291366577148

# --- Connector Id=116 "Conector: Angulado 115" --------------------------
$conn116 = Get-ShapeById $s.Shapes 116
$conn116.Left   = 405.4460601806641
$conn116.Top    = 54.66031455993652
$conn116.Width  = 77.04928970336915
$conn116.Height = 180.74362945556643
